$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.269.03"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.640.36"
$ws.Range("E3").Value = "  +2.64%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'594.30"
$ws.Range("E5").Value = "  +2.47%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'143.55"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.44%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.639.92"
$ws.Range("E9").Value = "  +2.71%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.24%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.67"
$ws.Range("E11").Value = "  +2.03%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.82%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +0.98%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'27.38"
$ws.Range("E14").Value = "  +2.27%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.115.45"
$ws.Range("E15").Value = "  +2.70%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "63.175.20"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +0.08%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.630.47"
$ws.Range("E18").Value = "  +2.10%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'11.37"
$ws.Range("E19").Value = "  +1.98%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'339.03"
$ws.Range("E20").Value = "  +0.50%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  +0.85%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'6.71"
$ws.Range("E22").Value = "  +1.19%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.05%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'66.85"
$ws.Range("E24").Value = "  -0.60%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  +5.70%  "

# Row 26 - SuiNetwork
$ws.Range("D26").Value = "'1.54"
$ws.Range("E26").Value = "  +2.74%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -0.16%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  +0.28%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'8.43"
$ws.Range("E29").Value = "  +2.76%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "'7.77"
$ws.Range("E30").Value = "  -1.87%  "

# Row 31 - Bittensor
$ws.Range("D31").Value = "'524.29"
$ws.Range("E31").Value = "  +14.68%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.82"
$ws.Range("E32").Value = "  +13.08%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +2.50%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0806"
$ws.Range("E34").Value = "  +0.56%  "

# Row 35 - Monero
$ws.Range("D35").Value = "'174.34"
$ws.Range("E35").Value = "  -1.11%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "'4.91"
$ws.Range("E36").Value = "  +11.12%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.02%  "

# Row 38 - PolygonEcosystemToken
$ws.Range("E38").Value = "  +1.65%  "

# Row 39 - EthereumClassic
$ws.Range("E39").Value = "  +0.98%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "'1.80"
$ws.Range("E40").Value = "  +7.67%  "

# Row 41 - Aave
$ws.Range("D41").Value = "'171.46"
$ws.Range("E41").Value = "  +8.08%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  -0.04%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'40.04"
$ws.Range("E43").Value = "  +0.14%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +1.47%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "'21.92"
$ws.Range("E45").Value = "  +5.37%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "'0.0557"
$ws.Range("E46").Value = "  +4.68%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "'0.629"
$ws.Range("E47").Value = "  +0.52%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +0.02%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +2.10%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'18.46"
$ws.Range("E50").Value = "  +2.84%  "

# Row 51 - WhiteBITCoin -> dogwifhat (coin replaced)
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.70"
$ws.Range("E51").Value = "  +1.95%  "
